# This edit swaps the data of row 2 and row 3 on the "Artfynd" sheet
# (all the per-observation fields that differ between the two species
# records), while leaving columns whose values happen to coincide
# between the two rows untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose content must be exchanged between row 2 and row 3.
$swapCols = @("A","B","E","F","G","H","P","Q","R","S","AW","AX")

foreach ($col in $swapCols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $val2 = $cell2.Value()
    $val3 = $cell3.Value()
    $cell2.Value = $val3
    $cell3.Value = $val2
}

# Columns Y and AA hold dates stored as plain text (e.g. "2020-03-26").
# Swap their text directly and force a text number format first so
# Excel does not reinterpret the strings as real date values.
$dateCols = @("Y","AA")
foreach ($col in $dateCols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $text2 = $cell2.Text
    $text3 = $cell3.Text
    $cell2.NumberFormat = "@"
    $cell3.NumberFormat = "@"
    $cell2.Value = $text3
    $cell3.Value = $text2
}
$ws.Range("Y2:Y3").ClearFormats()
$ws.Range("AA2:AA3").ClearFormats()

# Column M ("Aktivitet") only has a value on row 2 ("färska spår").
# After the edit it belongs to row 3 instead, and row 2's cell becomes
# empty again.
$activity = $ws.Range("M2").Value()
$ws.Range("M3").Value = $activity
$ws.Range("M2").ClearContents()
